$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("J7").Value = 2.82
$ws.Range("L7").Value = 3.65
$ws.Range("W7").Value = 6.8
$ws.Range("AA7").Value = 20
$ws.Range("AC7").Value = 8.25
$ws.Range("AD7").Value = 6
$ws.Range("AF7").Value = 70
$ws.Range("AG7").Value = 9.25
$ws.Range("AH7").Value = 17
$ws.Range("AK7").Value = 29
$ws.Range("AL7").Value = 37
$ws.Range("AM7").Value = 500
$ws.Range("AO7").Value = 11.5
$ws.Range("AP7").Value = 20
$ws.Range("AQ7").Value = 50
$ws.Range("AR7").Value = 80
$ws.Range("AW7").Value = 5
$ws.Range("AX7").Value = 17
$ws.Range("AY7").Value = 23
$ws.Range("AZ7").Value = 80
$ws.Range("BA7").Value = 110

# Row 10
$ws.Range("I10").Value = 3.3
$ws.Range("K10").Value = 1.91
$ws.Range("O10").Value = 1.53
$ws.Range("P10").Value = 2.38
$ws.Range("Q10").Value = 2.7
$ws.Range("R10").Value = 1.44
$ws.Range("X10").Value = 11
$ws.Range("AF10").Value = 67
$ws.Range("AS10").Value = 301
$ws.Range("BA10").Value = 101

# Row 11
$ws.Range("M11").Value = 1.08
$ws.Range("N11").Value = 8
$ws.Range("Q11").Value = 2.2
$ws.Range("R11").Value = 1.65

# Row 13
$ws.Range("I13").Value = 4.8
$ws.Range("J13").Value = 2.22
$ws.Range("L13").Value = 5
$ws.Range("P13").Value = 2.82
$ws.Range("W13").Value = 6.2
$ws.Range("X13").Value = 7.6
$ws.Range("AA13").Value = 14
$ws.Range("AB13").Value = 29
$ws.Range("AE13").Value = 17
$ws.Range("AF13").Value = 90
$ws.Range("AG13").Value = 11.5
$ws.Range("AH13").Value = 27
$ws.Range("AI13").Value = 16
$ws.Range("AK13").Value = 55
$ws.Range("AL13").Value = 60
$ws.Range("AM13").Value = 800
$ws.Range("AN13").Value = 3.5
$ws.Range("AO13").Value = 8
$ws.Range("AP13").Value = 17
$ws.Range("AQ13").Value = 27
$ws.Range("AR13").Value = 55
$ws.Range("AU13").Value = 7.4
$ws.Range("AV13").Value = 70
$ws.Range("AX13").Value = 28
$ws.Range("AY13").Value = 32
$ws.Range("AZ13").Value = 175
$ws.Range("BA13").Value = 200
$ws.Range("BB13").Value = 450

$wb.Save()
